$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 121 ---
# Was: A121 = "Print all nodes that are at distance k from a leaf node"
#      B121 = "https://gist.github.com/baranis/3f06f159aacc49697fac"
# Now: A121 = "UniqueNodes at K Distance From Root,Assuming No Duplicate Values at Nodes"
#      C121 = hyperlink -> http://ideone.com/ymyLOU  (B121 cleared out, link moved to column C)
$ws.Range("A121").Value = "UniqueNodes at K Distance From Root,Assuming No Duplicate Values at Nodes"
$ws.Range("B121").Clear()
$ws.Hyperlinks.Add($ws.Range("C121"), "http://ideone.com/ymyLOU") | Out-Null
# Hyperlinks.Add stamps the built-in "Hyperlink" style on the cell; realign it
# with the style already used by the sheet's other link cells (e.g. C115).
$ws.Range("C121").Style = $ws.Range("C115").Style

# --- Row 122 ---
# Was: A122 = "ReverseEveryKNodes"
#      B122 = "https://gist.github.com/baranis/1ba89234698fbe642e04"
# Now: A122 = "Reverse Every K Nodes in Linked List"
#      C122 = hyperlink -> http://ideone.com/25I7AF  (B122's hyperlink moved to column C)
$ws.Range("A122").Value = "Reverse Every K Nodes in Linked List"
$ws.Range("B122").Clear()
$ws.Hyperlinks.Add($ws.Range("C122"), "http://ideone.com/25I7AF") | Out-Null
$ws.Range("C122").Style = $ws.Range("C115").Style

# --- Row 123 removed entirely ---
# Was: A123 = "GetAllUniqueNodesKDistanceFromRoot"
#      B123 = "https://gist.github.com/baranis/b8d75d53c53734176d2e"
# That whole row/problem entry is dropped (merged into row 121's update above).
$ws.Rows.Item(123).Delete()
